$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08785
$ws.Range("H2").Value = 0.26355
$ws.Range("M2").Value = 2.931802
$ws.Range("N2").Value = 8.795406
$ws.Range("O2").Value = 0.2031783997257507
$ws.Range("P2").Value = 0.2031783997257507
$ws.Range("Q2").Value = 0.2575588056999999
$ws.Range("R2").Value = 2.3180292513
$ws.Range("S2").Value = 0.2031783997257507
$ws.Range("T2").Value = 0.2031783997257507

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08785
$ws.Range("H3").Value = 0.26355
$ws.Range("O3").Value = 0.1045038148188874
$ws.Range("P3").Value = 0.1045038148188874
$ws.Range("Q3").Value = 0.1324741103
$ws.Range("R3").Value = 1.1922669927
$ws.Range("S3").Value = 0.1045038148188874
$ws.Range("T3").Value = 0.1045038148188874

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08785
$ws.Range("H4").Value = 0.26355
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3072716666666667
$ws.Range("N4").Value = 0.9218150000000001
$ws.Range("O4").Value = 0.02129440034299643
$ws.Range("P4").Value = 0.02129440034299643
$ws.Range("Q4").Value = 0.02699381591666667
$ws.Range("R4").Value = 0.24294434325
$ws.Range("S4").Value = 0.02129440034299643
$ws.Range("T4").Value = 0.02129440034299643

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08785
$ws.Range("H5").Value = 0.26355
$ws.Range("M5").Value = 0.7995056666666667
$ws.Range("N5").Value = 2.398517
$ws.Range("O5").Value = 0.05540697561602141
$ws.Range("P5").Value = 0.05540697561602141
$ws.Range("Q5").Value = 0.07023657281666666
$ws.Range("R5").Value = 0.63212915535
$ws.Range("S5").Value = 0.05540697561602141
$ws.Range("T5").Value = 0.05540697561602141

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08785
$ws.Range("H6").Value = 0.26355
$ws.Range("M6").Value = 8.335887
$ws.Range("N6").Value = 25.007661
$ws.Range("O6").Value = 0.5776898238539604
$ws.Range("P6").Value = 0.5776898238539604
$ws.Range("Q6").Value = 0.73230767295
$ws.Range("R6").Value = 6.59076905655
$ws.Range("S6").Value = 0.5776898238539604
$ws.Range("T6").Value = 0.5776898238539604

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08785
$ws.Range("H7").Value = 0.26355
$ws.Range("M7").Value = 0.547269
$ws.Range("N7").Value = 1.641807
$ws.Range("O7").Value = 0.03792658564238371
$ws.Range("P7").Value = 0.03792658564238371
$ws.Range("Q7").Value = 0.04807758165
$ws.Range("R7").Value = 0.43269823485
$ws.Range("S7").Value = 0.03792658564238371
$ws.Range("T7").Value = 0.03792658564238371
